# Reorders / edits the Requisitos (course requirements) list in the Word
# document. All 21 "LOxxxx - ..." requirement lines live as separate runs
# (text + line <w:br/>) inside a single ListBullet-styled paragraph. The
# commit reorders them, drops "LOB1045 - Leitura e Producao de Textos
# Academicos", adds "LOB1012 - Estatistica", and fixes the "Algebra Linear"
# accent typo. We rebuild that paragraph's run content in one shot via
# Range.InsertXML so the resulting OOXML has one <w:r> per item, matching
# the target shape exactly.

$d = $word.ActiveDocument

# Locate the "Requisitos" heading paragraph, then the very next paragraph
# (the ListBullet paragraph holding all the requirement lines).
$reqHeadingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Requisitos") {
        $reqHeadingIndex = $i
    }
}
if ($reqHeadingIndex -eq -1) {
    throw "Could not find the 'Requisitos' heading paragraph"
}

$listPara = $d.Paragraphs.Item($reqHeadingIndex + 1)
$listRange = $listPara.Range
$start = $listRange.Start
$end = $listRange.End

$targetRange = $d.Range($start, $end)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>LOQ4095 -  Química Geral Experimental  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1006 -  Cálculo IV  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1040 -  Laboratório de Eletricidade  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1053 -  Física III  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1003 -  Cálculo I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1011 -  Eletricidade Aplicada  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1012 -  Estatística  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1024 -  Mecânica  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1036 -  Geometria Analítica  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1037 -  Álgebra Linear  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1038 -  Física Experimental I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1039 -  Física Experimental III  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1052 -  Cálculo III  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1004 -  Cálculo II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1018 -  Física I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1019 -  Física II  (Requisito)</w:t><w:br/></w:r></w:p>'

$targetRange.InsertXML($newParaXml)

Write-Output ("Updated Requisitos list; now " + $d.Paragraphs.Item($reqHeadingIndex + 1).Range.Text.Length + " chars")
